$d = $word.ActiveDocument
$d.Content.Find.Execute("wit hand without", $false, $false, $false, $false, $false, $true, 1, $false, "with and without", 2)
